$d = $word.ActiveDocument

# Locate the target paragraph ("ACIONISTAS PRESENTES: ...") that holds the
# shareholders Jinja loop.
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -clike "*shareholders*") {
        $p = $para
        break
    }
}

# --- 1) `item.person_type` -> `item.type` (both occurrences in this paragraph) ---
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Find.Execute("item.person_type", $true, $false, $false, $false, $false, $true, 0, $false, "item.type", 2)

# --- 2) Drop the duplicated "{{ item.name.text | upper }}, " run that sat
#        between "pessoa juridica de direito privado," and
#        "neste ato devidamente representada pelos seus diretores ", merging
#        those two text runs into one (with a single space joining them). ---
$r = $d.Range($p.Range.Start, $p.Range.End)
$found = $r.Find.Execute("{{ item.name.text  | upper }}, ", $true, $false, $false, $false, $false, $true, 0, $false, "", 1)
if ($found) {
    $ins = $d.Range($r.Start, $r.Start)
    $ins.InsertAfter(" ")
}

# --- 3) Drop the first "{% for diretor in item %} {{ item.name.text | upper }}, "
#        block entirely, and insert a single bold/yellow
#        "{{ item.name.text | upper }}, " run right before
#        "{% endfor %}{% elif item.type == 'individual' %}". ---
$r = $d.Range($p.Range.Start, $p.Range.End)
$found = $r.Find.Execute("{% for diretor in item %} {{ item.name.text | upper }}, ", $true, $false, $false, $false, $false, $true, 0, $false, "", 1)
if ($found) {
    $insStart = $r.Start
    $ins = $d.Range($insStart, $insStart)
    $ins.InsertAfter("{{ item.name.text | upper }}, ")
    $newRng = $d.Range($insStart, $insStart + 30)
    $newRng.Bold = 1
    $newRng.Font.HighlightColorIndex = 7
}
